$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.148.22"
$ws.Range("D2").Style = $ws.Range("B2").Style
$ws.Range("E2").Value = "  -0.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.900.25"
$ws.Range("D3").Style = $ws.Range("B3").Style
$ws.Range("E3").Value = "  -0.06%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = $ws.Range("B4").Style
$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.05"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "  +0.16%  "

$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3801"
$ws.Range("D8").Style = $ws.Range("B8").Style
$ws.Range("E8").Value = "  +0.83%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07287"
$ws.Range("D9").Style = $ws.Range("B9").Style
$ws.Range("E9").Value = "  +0.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.33"
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = "  +0.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9060"
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = "  +0.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08218"
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = "  -1.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.891.81"
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("E13").Value = "  -0.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.51"
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = "  +0.91%  "

$ws.Range("E15").Value = "  +1.66%  "

$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008656"
$ws.Range("D17").Style = $ws.Range("B17").Style
$ws.Range("E17").Value = "  +0.89%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.66"
$ws.Range("D18").Style = $ws.Range("B18").Style
$ws.Range("E18").Value = "  +1.22%  "

$ws.Range("E19").Value = "  +0.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.189.65"
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = "  -0.05%  "

$ws.Range("E21").Value = "  +1.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.126.26"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = "  -1.12%  "

$ws.Range("E23").Value = "  +1.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.468"
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = "  +0.80%  "

$ws.Range("E25").Value = "  +2.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "149.60"
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = "  +2.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.28"
$ws.Range("D27").Style = $ws.Range("B27").Style

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.736"
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").Value = "  -1.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.41"
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value = "  +0.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.824"
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value = "  +0.88%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.871"
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value = "  -0.87%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09232"
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = "  +0.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05047"
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7941"
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Value = "  -3.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.225"
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Value = "  -0.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.982"
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = "  +0.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.382"
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = "  +0.56%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.643"
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Value = "  +2.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5737"
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Value = "  +0.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01993"
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value = "  +1.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.082"
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = "  +0.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.041"
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = "  +1.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.621"
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = "  -0.49%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "116.32"
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = "  -1.59%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1517"
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = "  +0.52%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4899"
$ws.Range("D46").Style = $ws.Range("B46").Style
$ws.Range("E46").Value = "  +1.68%  "

$ws.Range("E47").Value = "  +0.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.12"
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = "  -0.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.641"
$ws.Range("D49").Style = $ws.Range("B49").Style

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "38.60"
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = "  +3.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.16"
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = "  +0.97%  "
